# Applies the scheduled-runner update to Golem_Profits sheets (currentAveragePrice /
# LevePrice / LeveProfit recompute) across ALC, ARM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1600.3334
$ws.Range("I137").Value = 901
$ws.Range("J137").Value = 1950
$ws.Range("K137").Value = 2703
$ws.Range("L137").Value = 5850
$ws.Range("M137").Value = -153
$ws.Range("N137").Value = -10950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1989
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1989
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1989
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2413

$ws.Range("H74").Value = 795
$ws.Range("I74").Value = 795
$ws.Range("K74").Value = 795
$ws.Range("M74").Value = 79

$ws.Range("H77").Value = 795
$ws.Range("I77").Value = 795
$ws.Range("K77").Value = 3975
$ws.Range("M77").Value = 393

$ws.Range("H136").Value = 1989
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1989
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 5967
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -11067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 127.5
$ws.Range("I4").Value = 127.5
$ws.Range("K4").Value = 127.5
$ws.Range("M4").Value = -15.5

$ws.Range("H22").Value = 820.4737
$ws.Range("I22").Value = 978.3333
$ws.Range("J22").Value = 549.8570999999999
$ws.Range("K22").Value = 978.3333
$ws.Range("L22").Value = 549.8570999999999
$ws.Range("M22").Value = -628.3333
$ws.Range("N22").Value = -1249.8571

$ws.Range("H31").Value = 1746.9
$ws.Range("I31").Value = 941.05554
$ws.Range("J31").Value = 8999.5
$ws.Range("K31").Value = 941.05554
$ws.Range("L31").Value = 8999.5
$ws.Range("M31").Value = -646.05554
$ws.Range("N31").Value = -9589.5

$ws.Range("H34").Value = 1746.9
$ws.Range("I34").Value = 941.05554
$ws.Range("J34").Value = 8999.5
$ws.Range("K34").Value = 941.05554
$ws.Range("L34").Value = 8999.5
$ws.Range("M34").Value = -739.05554
$ws.Range("N34").Value = -9403.5

$ws.Range("H35").Value = 5588
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20588

$ws.Range("H47").Value = 35328.332
$ws.Range("J47").Value = 35492.5
$ws.Range("L47").Value = 35492.5
$ws.Range("N47").Value = -36624.5

$ws.Range("H58").Value = 2000
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 2000
$ws.Range("N58").Value = -2406

$ws.Range("H59").Value = 74962.25
$ws.Range("J59").Value = 74962.25
$ws.Range("L59").Value = 74962.25
$ws.Range("N59").Value = -77252.25

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 2000
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8506.25
$ws.Range("I3").Value = 8008.3335
$ws.Range("K3").Value = 24025.0005
$ws.Range("M3").Value = -23913.0005

$ws.Range("H44").Value = 884.5625
$ws.Range("I44").Value = 241.4
$ws.Range("J44").Value = 1003.6667
$ws.Range("K44").Value = 724.2
$ws.Range("L44").Value = 3011.0001
$ws.Range("M44").Value = -326.2
$ws.Range("N44").Value = -3807.0001

$ws.Range("H55").Value = 3684.2205
$ws.Range("J55").Value = 3912.7358
$ws.Range("L55").Value = 11738.2074
$ws.Range("N55").Value = -12092.2074

$ws.Range("H64").Value = 1923.75

$ws.Range("H67").Value = 1923.75

$ws.Range("H108").Value = 353.25
$ws.Range("I108").Value = 353.25
$ws.Range("K108").Value = 1059.75
$ws.Range("M108").Value = 1820.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 143.28572
$ws.Range("I2").Value = 153.21053
$ws.Range("J2").Value = 122.333336
$ws.Range("K2").Value = 153.21053
$ws.Range("L2").Value = 122.333336
$ws.Range("M2").Value = -40.21053000000001
$ws.Range("N2").Value = -348.333336

$ws.Range("H47").Value = 29797.5
$ws.Range("J47").Value = 29797.5
$ws.Range("L47").Value = 29797.5
$ws.Range("N47").Value = -30933.5

$ws.Range("H80").Value = 5725
$ws.Range("I80").Value = 5725
$ws.Range("J80").Value = 5725
$ws.Range("K80").Value = 5725
$ws.Range("L80").Value = 5725
$ws.Range("M80").Value = -4727
$ws.Range("N80").Value = -7721

$ws.Range("H83").Value = 5725
$ws.Range("I83").Value = 5725
$ws.Range("J83").Value = 5725
$ws.Range("K83").Value = 28625
$ws.Range("L83").Value = 28625
$ws.Range("M83").Value = -23633
$ws.Range("N83").Value = -38609

$ws.Range("H126").Value = 4550.8
$ws.Range("I126").Value = 4253.6665
$ws.Range("J126").Value = 4996.5
$ws.Range("K126").Value = 12760.9995
$ws.Range("L126").Value = 14989.5
$ws.Range("M126").Value = -10290.9995
$ws.Range("N126").Value = -19929.5

$ws.Range("H132").Value = 2189.2856
$ws.Range("I132").Value = 1862.2
$ws.Range("K132").Value = 5586.6
$ws.Range("M132").Value = -3056.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1228.7142
$ws.Range("I32").Value = 1228.7142
$ws.Range("K32").Value = 1228.7142
$ws.Range("M32").Value = -911.7141999999999

$ws.Range("H68").Value = 4125
$ws.Range("I68").Value = 4666.6665
$ws.Range("K68").Value = 4666.6665
$ws.Range("M68").Value = -3917.6665

$ws.Range("H71").Value = 4125
$ws.Range("I71").Value = 4666.6665
$ws.Range("K71").Value = 23333.3325
$ws.Range("M71").Value = -19589.3325

$ws.Range("H82").Value = 3170.6667
$ws.Range("I82").Value = 2149.8572
$ws.Range("J82").Value = 4599.8
$ws.Range("K82").Value = 2149.8572
$ws.Range("L82").Value = 4599.8
$ws.Range("M82").Value = -1788.8572
$ws.Range("N82").Value = -5321.8

$ws.Range("H85").Value = 3170.6667
$ws.Range("I85").Value = 2149.8572
$ws.Range("J85").Value = 4599.8
$ws.Range("K85").Value = 2149.8572
$ws.Range("L85").Value = 4599.8
$ws.Range("M85").Value = -901.8571999999999
$ws.Range("N85").Value = -7095.8

$ws.Range("H93").Value = 3369.375
$ws.Range("I93").Value = 3680
$ws.Range("J93").Value = 2851.6667
$ws.Range("K93").Value = 3680
$ws.Range("L93").Value = 2851.6667
$ws.Range("M93").Value = -2432
$ws.Range("N93").Value = -5347.6667

$ws.Range("H132").Value = 8333
$ws.Range("I132").Value = 8333
$ws.Range("K132").Value = 24999
$ws.Range("M132").Value = -22469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 500500
$ws.Range("I2").Value = 500500
$ws.Range("K2").Value = 500500
$ws.Range("M2").Value = -500388

$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 4000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -4226

$ws.Range("H69").Value = 10166
$ws.Range("J69").Value = 10166
$ws.Range("L69").Value = 10166
$ws.Range("N69").Value = -11664

$ws.Range("H72").Value = 10166
$ws.Range("J72").Value = 10166
$ws.Range("L72").Value = 30498
$ws.Range("N72").Value = -37986
